# Add an "Electrode Locations" column (C) to the main dataframe, and
# sort all data rows by electrode location (A1 .. O15 -> letter then
# numeric order), mirroring the source commit's change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine how many data rows currently exist (header is row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Collect the existing File Name / Unnormalized P_max pairs together with
# a derived "electrode location" (e.g. "A11_bipolar_10V_1kHz.txt" -> "A11")
# and a zero-padded sort key so numeric suffixes sort naturally
# (A2 < A11, not "A11" < "A2" as plain text would).
$records = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $fileName = $ws.Cells.Item($r, 1).Text
    $pVal = $ws.Cells.Item($r, 2).Value2

    if ($fileName -match '^([A-Za-z]+)(\d+)_') {
        $letter = $matches[1]
        $number = [int]$matches[2]
    } else {
        $letter = $fileName
        $number = 0
    }

    $location = "$letter$number"
    $sortKey = "{0}_{1:D4}" -f $letter, $number

    $records += [PSCustomObject]@{
        FileName = $fileName
        PMax     = $pVal
        Location = $location
        SortKey  = $sortKey
    }
}

$sortedRecords = $records | Sort-Object -Property SortKey

# Header for the new column, matching the header style already used by
# "File Name" / "Unnormalized P_max" (A1).
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Electrode Locations"

# Re-write columns A, B and C in the newly sorted order.
$row = 2
foreach ($rec in $sortedRecords) {
    $ws.Cells.Item($row, 1).Value = $rec.FileName
    $ws.Cells.Item($row, 2).Value = $rec.PMax
    $ws.Cells.Item($row, 3).Value = $rec.Location
    $row++
}
